# Apply "add the final report format" edit:
#  - Delete rows 5 and 6 (the two extra "Items to be improved" entries)
#  - Populate column B (Description) for the remaining data rows
#  - Rewrite column C (Improvement direction) text for the remaining rows
#  - Refresh column G (Expected completion date) wording for row 3
#  - Refresh column J (Confirmation) wording for rows 3 and 4
#  - Fill in column K (appendix) for row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing item rows (old rows 5 & 6) first so the remaining
# rows keep their original row numbers (2, 3, 4).
$ws.Rows("5:6").Delete()

# Row 2 - "The saturation of the lower glass point Xiaoli Pill is 65%"
$ws.Range("B2").Value = "Low job saturation(lower than95%)"
$ws.Range("C2").Value = "1. Judgment of value and reduce tasks without added value`n2. Inspection of movement quality and human engineering hazards: Reduce the number of movements, work with both hands at the same time, shorten the distance of movements, and make movements easier; eliminate human engineering hazards`n3. Automated level inspection: simple and automated import`n4. Merge and rearrange new job elements"
$ws.Range("K2").Value = "LiXX"

# Row 3 - "Bottom glass electrophoresis tank+UVFixed baking operation saturation79.2%"
$ws.Range("B3").Value = "Low job saturation(lower than95%)"
$ws.Range("C3").Value = "The improvement direction for the problem 'Bottom glass electrophoresis tank+UVFixed baking operation saturation79.2%' is as follows: 1. Judgment of value and reduce tasks without added value 2. Inspection of movement quality and human engineering hazards: Reduce the number of movements, work with both hands at the same time, shorten the distance of movements, and make movements easier; eliminate human engineering hazards 3. Automated level inspection: simple and automated import 4. Merge and rearrange new job elements."
$ws.Range("G3").Value = "Expected completion date: 10/15/24"
$ws.Range("J3").Value = "The improvement parameter 'Confirmation' for the problem 'Bottom glass electrophoresis tank+UVFixed baking operation saturation79.2%' is 'yes'."

# Row 4 - "Xiaoliwan wax+Paste conductive foam+Lower glass glue frame dispensing operation saturation75%"
$ws.Range("B4").Value = "Low job saturation(lower than95%)"
$ws.Range("C4").Value = "1. Judgment of value and reduce tasks without added value`n2. Inspection of movement quality and human engineering hazards`n3. Automated level inspection`n4. Merge and rearrange new job elements"
$ws.Range("J4").Value = "The improvement parameter 'Confirmation' for the problem 'Xiaoliwan wax+Paste conductive foam+Lower glass glue frame dispensing operation saturation75%' is 'yes'."
